# Tilføjet bilag nr. til titlen på bilagene
# Prepend "Bilag 9 - " (as three separate runs) to the first paragraph's
# title text ("Fully-dressed Use Case") and bump the title's font size
# from 16pt (sz 32) to 20pt (sz 40).

$d = $word.ActiveDocument

# Insert the three new runs at the very start of the document, one at a
# time (each InsertBefore on a fresh zero-length range produces its own
# <w:r> run instead of merging into the neighbouring run).
$d.Range(0, 0).InsertBefore(" - ")
$d.Range(0, 0).InsertBefore("9")
$d.Range(0, 0).InsertBefore("Bilag ")

# Re-apply bold + the new 20pt (half-point 40) size across the whole
# title paragraph (new text + the pre-existing "Fully-dressed Use Case"
# text) so that both the ASCII (sz/szCs) and the paragraph mark's rPr
# all end up at size 40, matching the rest of the title's formatting.
$titleRange = $d.Paragraphs(1).Range
$titleRange.Font.Bold = $true
$titleRange.Font.BoldBi = $true
$titleRange.Font.Size = 20
$titleRange.Font.SizeBi = 20
